$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '26.538.27'; E = '  +7.02%  ' },
    @{ Row = 3; D = '1.742.30'; E = '  +5.27%  ' },
    @{ Row = 4; D = '1.006'; E = '  +0.13%  ' },
    @{ Row = 5; D = '334.71'; E = '  +7.71%  ' },
    @{ Row = 6; D = '1.003'; E = '  +0.17%  ' },
    @{ Row = 7; D = '0.3787'; E = '  +4.63%  ' },
    @{ Row = 8; D = '48.65'; E = '  +3.35%  ' },
    @{ Row = 9; D = '0.3385'; E = '  +4.50%  ' },
    @{ Row = 10; D = '1.186'; E = '  +5.42%  ' },
    @{ Row = 11; D = '0.07475'; E = '  +6.27%  ' },
    @{ Row = 12; D = '1.003'; E = '  +0.05%  ' },
    @{ Row = 13; D = '6.459'; E = '  +7.46%  ' },
    @{ Row = 14; D = '20.40'; E = '  +5.24%  ' },
    @{ Row = 15; D = '7.118'; E = '  +8.47%  ' },
    @{ Row = 16; D = '1.744.30'; E = '  +5.12%  ' },
    @{ Row = 17; D = '0.00001087'; E = '  +4.37%  ' },
    @{ Row = 18; D = '0.06710'; E = '  +2.03%  ' },
    @{ Row = 19; D = '83.45'; E = '  +6.18%  ' },
    @{ Row = 20; D = '1.002'; E = '  +0.07%  ' },
    @{ Row = 21; D = '16.73'; E = '  +6.83%  ' },
    @{ Row = 22; D = '6.208'; E = '  +5.61%  ' },
    @{ Row = 23; D = '13.05'; E = '  +4.67%  ' },
    @{ Row = 24; D = '26.564.81'; E = '  +7.07%  ' },
    @{ Row = 25; D = '2.452'; E = '  +0.89%  ' },
    @{ Row = 26; D = '2.466'; E = '  +1.09%  ' },
    @{ Row = 27; D = '1.413'; E = '  +18.71%  ' },
    @{ Row = 28; D = '154.33'; E = '  +4.80%  ' },
    @{ Row = 29; D = '19.67'; E = '  +6.27%  ' },
    @{ Row = 30; D = '1.935.76'; E = '  +5.14%  ' },
    @{ Row = 31; D = '132.37'; E = '  +5.97%  ' },
    @{ Row = 32; D = '4.142'; E = '  +1.80%  ' },
    @{ Row = 33; D = '6.105'; E = '  +6.79%  ' },
    @{ Row = 34; D = '0.08704'; E = $null },
    @{ Row = 35; D = '1.714'; E = '  +3.98%  ' },
    @{ Row = 36; D = '12.97'; E = '  +6.53%  ' },
    @{ Row = 37; D = '5.436'; E = '  +5.53%  ' },
    @{ Row = 38; D = '0.02355'; E = '  +4.85%  ' },
    @{ Row = 39; D = '0.06319'; E = '  +4.72%  ' },
    @{ Row = 40; D = '0.2180'; E = '  +5.87%  ' },
    @{ Row = 41; D = '8.595'; E = '  +3.58%  ' },
    @{ Row = 42; D = '1.230'; E = '  -4.05%  ' },
    @{ Row = 43; D = '0.6235'; E = '  +5.56%  ' },
    @{ Row = 44; D = '14.33'; E = '  +12.64%  ' },
    @{ Row = 45; D = $null; E = '  +0.18%  ' },
    @{ Row = 46; D = '3.930'; E = '  +4.54%  ' },
    @{ Row = 47; D = '0.6066'; E = '  +8.45%  ' },
    @{ Row = 48; D = '128.53'; E = '  +2.92%  ' },
    @{ Row = 49; D = '2.063'; E = '  +6.61%  ' },
    @{ Row = 50; D = '0.07246'; E = '  +4.12%  ' },
    @{ Row = 51; D = '77.88'; E = '  +4.50%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $cell = $ws.Range("E" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}